$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows right before the current row 290 (pushes the
# existing rows 290-305 down to 293-308, matching the XML diff's row
# renumbering) and carries formatting (incl. the date style on column D)
# down from the row above, same as native Excel row-insert behaviour.
$ws.Rows("290:292").Insert()

# New weekly price observations for Vega Monumental Concepción - Limón,
# fecha 44516 ("2021-11-16"), provincia de Melipilla, malla 16 kilos.
$newRows = @(
    @{ Row=290; Calidad="1a amarillo"; Volumen=300; PMin=7500; PMax=7500; PProm=7500; Origen="Provincia de Melipilla"; PrecioKg=469 },
    @{ Row=291; Calidad="1a plateado"; Volumen=300; PMin=8500; PMax=8500; PProm=8500; Origen="Provincia de Melipilla"; PrecioKg=531 },
    @{ Row=292; Calidad="2a amarillo"; Volumen=300; PMin=6500; PMax=6500; PProm=6500; Origen="Provincia de Melipilla"; PrecioKg=406 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = 44516
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = '$/malla 16 kilos'
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = 16
}
